$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column cell format (style s="2") down into the new row 19
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 4.930115226412357
$ws.Range("D2").Value = 2008

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 1.457587285166628
$ws.Range("D3").Value = 2009

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -0.9140166223623458
$ws.Range("D4").Value = 2010

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 2.585942866987878
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 4.109775046142405

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 4.253963781362402
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 2.863367440851095

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 1.752870900283909
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 1.520397254708405

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = -1.479696720105139
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 2.503951807923066

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 3.900127535411246
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 1.194160460927884

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 0.03947433952959933
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 1.459149667419779

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 2.192778679161944
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 1.586376095629216

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 3.40836448860673
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 2.570658574505469

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 2.799070570134488
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 2.479713128614147

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 4.195393191694419
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 2.359935293525561

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = 1.666553973046048
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = -0.4512719783814068

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 1.879266440112803
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 1.081814991510499

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = -2.620683231370946
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = -1.174318230871441

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -3.036556262700274
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.07123445333143685

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -2.953443685011514
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = -1.196842846539037
